$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-05 12:55:06"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
